$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing hyperlink target URLs (F2:F23) in row order before
# the new row shifts everything down. These will simply move down one row. ---
$oldUrls = @()
for ($r = 2; $r -le 23; $r++) {
    $oldUrls += $ws.Range("F" + $r).Value2
}

# New circular published 2025-12-24, becomes the new first data row (Sl.no. 23)
$newUrl = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-24-12-2025.pdf"

# --- Insert a new row above the current first data row (row 2), pushing all
# existing data rows (and their formatting) down by one. ---
$ws.Rows(2).Insert()

# Copy the format of the row right below (the old row-2 format, now at row 3)
# into the freshly inserted, blank row 2 so styles match the rest of the table.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Populate the new row 2 with the latest circular's data. ---
$ws.Range("A2").Value = 23
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 301.65
$ws.Range("E2").Value = "24-12-2025"
$ws.Range("F2").Value = $newUrl

# --- Rebuild the hyperlinks so the Circular Link column (F) keeps working
# after the row shift (the engine's row-insert does not itself relocate the
# <hyperlinks> entries), then add the new link for row 2. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), $newUrl) | Out-Null
for ($i = 0; $i -lt $oldUrls.Count; $i++) {
    $row = 3 + $i
    $ws.Hyperlinks.Add($ws.Range("F" + $row), $oldUrls[$i]) | Out-Null
}

# Hyperlinks.Add() stamps the built-in "Hyperlink" cell style (blue/underline)
# on its target cell; restore the plain centered look shared by the rest of
# the table by pasting column E's format (style index 3) over column F.
$ws.Range("E2:E24").Copy() | Out-Null
$ws.Range("F2:F24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
